# Apply the edits described by the commit diff to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Narrow the "_FilterDatabase" defined name from column V to column E
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='ExportedData (2)'!`$A`$1:`$E`$67"
    }
}

# ---------------------------------------------------------------------
# 2) Update the data cells in rows 65-67
# ---------------------------------------------------------------------
# Row 65: longitude value corrected
$ws.Range("B65").Value = 39.628986500000003

# Row 66: longitude value corrected, team lookup now fails (#N/A)
$ws.Range("B66").Value = 39.628986500000003
$ws.Range("D66").Value = "#N/A"

# Row 67: latitude/longitude populated, team lookup now fails (#N/A)
$ws.Range("A67").Value = 24.439459803763501
$ws.Range("A67").Style = "Normal"
$ws.Range("B67").Value = 39.628986500000003
$ws.Range("D67").Value = "#N/A"

# ---------------------------------------------------------------------
# 3) Shrink the AutoFilter range from column V to column E
# ---------------------------------------------------------------------
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("B1:E67").AutoFilter() | Out-Null

# ---------------------------------------------------------------------
# 4) Update the active selection shown when the sheet is opened
# ---------------------------------------------------------------------
$ws.Range("D13").Select() | Out-Null
